$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text format so numeric-looking
# strings (e.g. "1.00", "2.90") are not auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = '37.083.48'
$ws.Range("E2").Value = '  +0.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = '2.045.60'
$ws.Range("E3").Value = '  -3.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = '247.53'
$ws.Range("E5").Value = '  -3.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = '0.651'
$ws.Range("E6").Value = '  -2.84%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = '55.43'
$ws.Range("E8").Value = '  +16.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = '61.72'
$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.375'
$ws.Range("E10").Value = '  +0.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0753'
$ws.Range("E11").Value = '  +1.82%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = '14.87'
$ws.Range("E13").Value = '  +3.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = '2.343.13'
$ws.Range("E14").Value = '  -3.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.811'
$ws.Range("E15").Value = '  -4.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = '5.18'
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = '2.046.19'
$ws.Range("E17").Value = '  -3.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = '36.975.16'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = '71.78'
$ws.Range("E19").Value = '  -2.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0883'
$ws.Range("E20").Value = '  +5.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = '14.15'
$ws.Range("E21").Value = '  +5.22%  '

$ws.Range("B22").NumberFormat = "@"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '236.41'
$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("B23").NumberFormat = "@"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.23'
$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.12%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = '169.03'
$ws.Range("E26").Value = '  -1.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = '9.03'
$ws.Range("E27").Value = '  -2.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D28").Value = '20.01'
$ws.Range("E28").Value = '  -7.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("E29").Value = '  -3.28%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value = '4.57'
$ws.Range("E31").Value = '  +0.85%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +12.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("E33").Value = '  +2.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D34").Value = '4.31'
$ws.Range("E34").Value = '  +2.78%  '

$ws.Range("B35").NumberFormat = "@"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.24%  '

$ws.Range("B36").NumberFormat = "@"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.0879'
$ws.Range("E36").Value = '  -9.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D37").Value = '2.26'
$ws.Range("E37").Value = '  -4.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("E38").Value = '  -5.78%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.33'
$ws.Range("E39").Value = '  -2.13%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("B40").Value = 'Gas'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D40").Value = '16.58'
$ws.Range("E40").Value = '  -35.50%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").Value = '0.103'
$ws.Range("E41").Value = '  +22.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = '18.12'
$ws.Range("E42").Value = '  +11.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0222'
$ws.Range("E43").Value = '  -1.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = '1.13'
$ws.Range("E44").Value = '  -5.19%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '95.62'
$ws.Range("E45").Value = '  -3.60%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").Value = '2.78'
$ws.Range("E46").Value = '  -1.18%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.16'
$ws.Range("E47").Value = '  +53.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D48").Value = '1.293.11'
$ws.Range("E48").Value = '  -5.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("E49").Value = '  +2.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("E50").Value = '  +1.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("D51").Value = '6.74'
$ws.Range("E51").Value = '  -5.48%  '
